# Figure 2, panel-tag shuffle: the labels "(b)" and "(c)" swap places.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shapeB = $null
$shapeC = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $txt = $shape.TextFrame.TextRange.Text
        if ($txt -eq "(b)") {
            $shapeB = $shape
        } elseif ($txt -eq "(c)") {
            $shapeC = $shape
        }
    }
}

if ($shapeB -ne $null) { $shapeB.TextFrame.TextRange.Text = "(c)" }
if ($shapeC -ne $null) { $shapeC.TextFrame.TextRange.Text = "(b)" }
